$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate save_data to use K (strikeouts) instead of Strike# as the
# source for column G, recalculated from the std/mean and s_vals pass.
$ws.Range("G2").Value = 0
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 4
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
